# updated statbar xpaths & diagnosis testcases
#
# Adds three new sheets to the workbook:
#   - CypherOutput_Message : a duplicate of the existing "Message" sheet
#   - StatOutput            : the stat-bar counts (files/sample/cases/study)
#   - StatOutput_Message    : the connection/cypher "message" block, repeated
#                              twice, where the 2nd Cypher line is the new
#                              stat-bar query text
#
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Re-usable text blocks (shared with the existing "Message" sheet, which
# is why Excel will fold these back onto the very same shared-string
# entries already present in the workbook).
# ---------------------------------------------------------------------
$neo4jUrlLabel = "Neo4j_URL:"
$neo4jUrlValue = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userLabel     = "User_name:"
$userValue     = "neo4j"
$pwdLabel      = "PWD:"
$pwdValue      = "icdcDBneo4j0"
$cypherLabel   = "Cypher:"
$oldCypher     = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Multicentric lymphoma''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$newCypher     = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Multicentric lymphoma'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$outputLabel   = "Output:"
$outputValue   = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC14_Canine_Filter_Diagnosis-MultiLymph_Neo4jData.xlsx"

function Fill-MessageSheet($ws, $cypherText) {
    $ws.Range("A1").Value = $neo4jUrlLabel
    $ws.Range("A2").Value = $neo4jUrlValue
    $ws.Range("A3").Value = $userLabel
    $ws.Range("A4").Value = $userValue
    $ws.Range("A5").Value = $pwdLabel
    $ws.Range("A6").Value = $pwdValue
    $ws.Range("A7").Value = $cypherLabel
    $ws.Range("A8").Value = $cypherText
    $ws.Range("A9").Value = $outputLabel
    $ws.Range("A10").Value = $outputValue
}

# ---------------------------------------------------------------------
# 1) CypherOutput_Message  -- duplicate of "Message"
# ---------------------------------------------------------------------
$cypherOutputMessage = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$cypherOutputMessage.Name = "CypherOutput_Message"
Fill-MessageSheet $cypherOutputMessage $oldCypher

# ---------------------------------------------------------------------
# 2) StatOutput  -- counts returned by the stat-bar cypher query
# ---------------------------------------------------------------------
$statOutput = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$statOutput.Name = "StatOutput"
$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

foreach ($addr in @("A2", "B2")) {
    $cell = $statOutput.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "0"
    $cell.ClearFormats()
}
foreach ($addr in @("C2", "D2")) {
    $cell = $statOutput.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "1"
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 3) StatOutput_Message  -- connection block repeated twice, 2nd Cypher
#    line updated to the new stat-bar query
# ---------------------------------------------------------------------
$statOutputMessage = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$statOutputMessage.Name = "StatOutput_Message"
Fill-MessageSheet $statOutputMessage $oldCypher

$statOutputMessage.Range("A11").Value = $neo4jUrlLabel
$statOutputMessage.Range("A12").Value = $neo4jUrlValue
$statOutputMessage.Range("A13").Value = $userLabel
$statOutputMessage.Range("A14").Value = $userValue
$statOutputMessage.Range("A15").Value = $pwdLabel
$statOutputMessage.Range("A16").Value = $pwdValue
$statOutputMessage.Range("A17").Value = $cypherLabel
$statOutputMessage.Range("A18").Value = $newCypher
$statOutputMessage.Range("A19").Value = $outputLabel
$statOutputMessage.Range("A20").Value = $outputValue

# ---------------------------------------------------------------------
# Keep the original first tab ("CypherOutput") as the selected/active
# sheet, same as before the edit.
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Select()
